# Update "想去人数" (F) and "最低票价" (G) figures on both the "展览" sheet
# and the "全部类型" sheet, which contain duplicated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value  = 1791
    $ws.Range("F4").Value  = 34
    $ws.Range("F5").Value  = 773
    $ws.Range("F7").Value  = 107
    $ws.Range("F8").Value  = 34
    $ws.Range("F9").Value  = 107
    $ws.Range("F13").Value = 126
    $ws.Range("F15").Value = 4174
    $ws.Range("F18").Value = 461
    $ws.Range("F19").Value = 396
    $ws.Range("F20").Value = 968
    $ws.Range("F21").Value = 1312
    $ws.Range("F22").Value = 360
    $ws.Range("F23").Value = 36
    $ws.Range("F24").Value = 41
    $ws.Range("F25").Value = 1926
    $ws.Range("F26").Value = 62
    $ws.Range("G26").Value = 45
    $ws.Range("F27").Value = 54
    $ws.Range("F28").Value = 87
    $ws.Range("F29").Value = 191
}
